# "Fruta / hortaliza, semanal"
# A new weekly observation was inserted into the "Ajo" (Vega Monumental
# Concepción) log as row 276, pushing every existing record from row 276
# down through row 331 one row further (new rows 277-332). The sheet's
# used range grows from A1:R331 to A1:R332.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 276; Excel shifts rows 276:331 down to 277:332
# and extends the used range to A1:R332 automatically.
$ws.Rows.Item(276).Insert()

# Populate the newly-inserted row 276 with the new record's data.
$ws.Cells.Item(276, 1).Value = 11
$ws.Cells.Item(276, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(276, 3).Value = "Bíobío"
$ws.Cells.Item(276, 4).Value = 45204
$ws.Cells.Item(276, 5).Value = 8
$ws.Cells.Item(276, 6).Value = 100112003
$ws.Cells.Item(276, 7).Value = "Ajo"
$ws.Cells.Item(276, 8).Value = "Chino"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 300
$ws.Cells.Item(276, 11).Value = 21000
$ws.Cells.Item(276, 12).Value = 22000
$ws.Cells.Item(276, 13).Value = 21667
$ws.Cells.Item(276, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(276, 15).Value = "China"
$ws.Cells.Item(276, 16).Value = 2167
$ws.Cells.Item(276, 17).Value = 10
$ws.Cells.Item(276, 18).Value = "Hortaliza"
